# "added in fixed effects"
#
# A new worksheet "zip2" is inserted between the existing "zip" and "county"
# sheets. It holds a new regression table (zip-level model with incrementally
# added fixed effects) that replaces the earlier intercept-based "zip" table.
# The "county" sheet itself is untouched content-wise; it's simply no longer
# the active/selected tab.

$wb = $excel.ActiveWorkbook

$zipSheet    = $wb.Worksheets.Item("zip")
$countySheet = $wb.Worksheets.Item("county")

# Capture/restore county's own selection state before doing anything else,
# since copying another sheet afterwards can clobber whichever sheet
# currently holds the selection.
$countySheet.Range("D30").Select()

# Duplicate the "zip" sheet; Excel places the copy immediately after it and
# hands out the next unused sheetId (3), which matches the target workbook.
$zipSheet.Copy($null, $zipSheet)

$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "zip2"

# Remove everything carried over from the "zip" sheet copy.
$newSheet.Cells.Clear()

# Populate "zip2" with the new fixed-effects regression table.
$arr = New-Object "object[,]" 23,6

$arr[0,0]  = "Regressor";                                 $arr[0,1]  = "Model 1";   $arr[0,2]  = "Model 2";   $arr[0,3]  = "Model 3";   $arr[0,4]  = "Model 4";   $arr[0,5]  = "Model 5"
$arr[1,0]  = "minority_percent";                          $arr[1,1]  = "0.086***";  $arr[1,2]  = "0.065***";  $arr[1,3]  = "0.048***";  $arr[1,4]  = "0.022***";  $arr[1,5]  = "0.021"
$arr[2,0]  = " ";                                         $arr[2,1]  = "(0.006)";   $arr[2,2]  = "(0.004)";   $arr[2,3]  = "(0.003)";   $arr[2,4]  = "(0.003)";   $arr[2,5]  = "(0.020)"
$arr[3,0]  = "minority_percent^2";                        $arr[3,1]  = " ";         $arr[3,2]  = " ";         $arr[3,3]  = " ";         $arr[3,4]  = " ";         $arr[3,5]  = "0.000*"
$arr[4,0]  = " ";                                         $arr[4,1]  = " ";         $arr[4,2]  = " ";         $arr[4,3]  = " ";         $arr[4,4]  = " ";         $arr[4,5]  = "(0.000)"
$arr[5,0]  = "minority_percent * median_family_income";   $arr[5,1]  = " ";         $arr[5,2]  = " ";         $arr[5,3]  = " ";         $arr[5,4]  = " ";         $arr[5,5]  = "-0.000"
$arr[6,0]  = " ";                                         $arr[6,1]  = " ";         $arr[6,2]  = " ";         $arr[6,3]  = " ";         $arr[6,4]  = " ";         $arr[6,5]  = "(0.000)"
$arr[7,0]  = "minority_percent * cook_pvi";               $arr[7,1]  = " ";         $arr[7,2]  = " ";         $arr[7,3]  = " ";         $arr[7,4]  = " ";         $arr[7,5]  = "-0.000"
$arr[8,0]  = " ";                                         $arr[8,1]  = " ";         $arr[8,2]  = " ";         $arr[8,3]  = " ";         $arr[8,4]  = " ";         $arr[8,5]  = "(0.000)"
$arr[9,0]  = "minority_percent * rural";                  $arr[9,1]  = " ";         $arr[9,2]  = " ";         $arr[9,3]  = " ";         $arr[9,4]  = " ";         $arr[9,5]  = "0.011"
$arr[10,0] = " ";                                         $arr[10,1] = " ";         $arr[10,2] = " ";         $arr[10,3] = " ";         $arr[10,4] = " ";         $arr[10,5] = "(0.011)"
$arr[11,0] = "minority_percent * estimate_gini_index";    $arr[11,1] = " ";         $arr[11,2] = " ";         $arr[11,3] = " ";         $arr[11,4] = " ";         $arr[11,5] = "-0.019"
$arr[12,0] = " ";                                         $arr[12,1] = " ";         $arr[12,2] = " ";         $arr[12,3] = " ";         $arr[12,4] = " ";         $arr[12,5] = "(0.040)"
$arr[13,0] = "minority_percent * violent_crime_rate";     $arr[13,1] = " ";         $arr[13,2] = " ";         $arr[13,3] = " ";         $arr[13,4] = " ";         $arr[13,5] = "-0.000"
$arr[14,0] = " ";                                         $arr[14,1] = " ";         $arr[14,2] = " ";         $arr[14,3] = " ";         $arr[14,4] = " ";         $arr[14,5] = "(0.000)"
$arr[15,0] = "Demographic Variables";                     $arr[15,1] = "X";         $arr[15,2] = "X";         $arr[15,3] = "X";         $arr[15,4] = "X";         $arr[15,5] = "X"
$arr[16,0] = "Loan-specific Variables";                   $arr[16,1] = $null;       $arr[16,2] = "X";         $arr[16,3] = "X";         $arr[16,4] = "X";         $arr[16,5] = "X"
$arr[17,0] = "Economic Variables";                        $arr[17,1] = $null;       $arr[17,2] = $null;       $arr[17,3] = "X";         $arr[17,4] = "X";         $arr[17,5] = "X"
$arr[18,0] = "COVID-19 Variables";                        $arr[18,1] = $null;       $arr[18,2] = $null;       $arr[18,3] = $null;       $arr[18,4] = "X";         $arr[18,5] = "X"
$arr[19,0] = $null;                                       $arr[19,1] = $null;       $arr[19,2] = $null;       $arr[19,3] = $null;       $arr[19,4] = $null;       $arr[19,5] = $null
$arr[20,0] = "Num. obs.";                                 $arr[20,1] = "2061489";   $arr[20,2] = "2061489";   $arr[20,3] = "2061489";   $arr[20,4] = "2061489";   $arr[20,5] = "2061489"
$arr[21,0] = "R2 (full model)";                           $arr[21,1] = "0.050";     $arr[21,2] = "0.311";     $arr[21,3] = "0.556";     $arr[21,4] = "0.811";     $arr[21,5] = "0.811"
$arr[22,0] = "Adj. R2 (full model)";                      $arr[22,1] = "0.049";     $arr[22,2] = "0.310";     $arr[22,3] = "0.556";     $arr[22,4] = "0.811";     $arr[22,5] = "0.811"

$newSheet.Range("A1:F23").Value = $arr

# Match the target workbook's view state: "zip2" becomes the active/selected
# tab, with A2 selected.
$newSheet.Range("A2").Select()
$newSheet.Activate()
